$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row (72) for 11/12/2025, appended after the last existing data row (71).
# Force column A to be stored as plain text (matching the existing Date column
# cells) rather than letting Excel auto-convert the date-like string into a
# serial date number; then reset the style back to Normal so no extra
# number-format/style is left on the cell.
$row = 72
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "11/12/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.2016542154360835
$ws.Cells.Item($row, 3).Value = 0.7983457845639165
